$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row (row 4) for the "707Air" / "nueva empresa" company record.
# Columns: A=Razon Social, B=Usuario, C=NIT, D=Contraseña, E=E-Mail, F=Descripcion, G=Id
$ws.Range("A4").Value = "nueva empresa"
$ws.Range("B4").Value = "esto"
$ws.Range("E4").Value = "over@over.com"
$ws.Range("F4").Value = "descripcion"
$ws.Range("G4").Value = 3

# C4 ("6454654654") and D4 ("12345678") must stay TEXT even though they look like
# numbers. Assigning a numeric-looking string straight to .Value auto-converts it
# to a number, so build it as a text formula result first and paste the computed
# value back in - that keeps the cell's stored type as text without touching any
# cell styles/number formats.
$ws.Range("Z1").Formula = "=""6454654654"""
$ws.Range("Z1").Copy()
$ws.Range("C4").PasteSpecial(-4163)

$ws.Range("Z2").Formula = "=""12345678"""
$ws.Range("Z2").Copy()
$ws.Range("D4").PasteSpecial(-4163)

$ws.Range("Z1:Z2").Clear()
